$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.940.19"
$ws.Range("E2").Value = "  -2.36%  "
$ws.Range("D3").Value = "1.796.40"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.08"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5306"
$ws.Range("E7").Value = "  -1.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3867"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07450"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.44"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("E11").Value = "  -2.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.178"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.442"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.37"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "1.792.12"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.37"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06543"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9991"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.26"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.958"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "27.976.17"
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.089"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.78"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.14"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "2.005.99"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.302"
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.90"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.099"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1087"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.668"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.504"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06938"
$ws.Range("E35").Value = "  +6.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2207"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02276"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.392"
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.24"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.192"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6114"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.414"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.29"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.676"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5709"
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "124.44"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.177"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.914"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06807"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000294"
$ws.Range("E51").Value = "  +32.44%  "
